$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.392.92"
$ws.Range("E2").Value = "  +4.93%  "
$ws.Range("D3").Value = "1.814.61"
$ws.Range("E3").Value = "  +5.68%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "344.89"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "0.3808"
$ws.Range("E7").Value = "  +3.21%  "
$ws.Range("D8").Value = "0.3493"
$ws.Range("E8").Value = "  +4.53%  "
$ws.Range("D9").Value = "48.73"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").Value = "1.233"
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("D11").Value = "0.07714"
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "22.07"
$ws.Range("E13").Value = "  +9.87%  "
$ws.Range("D14").Value = "6.624"
$ws.Range("E14").Value = "  +5.19%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.817.09"
$ws.Range("E15").Value = "  +6.06%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.219"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").Value = "0.00001117"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("D18").Value = "0.06742"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").Value = "86.04"
$ws.Range("E19").Value = "  +5.10%  "
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "17.62"
$ws.Range("E21").Value = "  +7.48%  "
$ws.Range("D22").Value = "6.586"
$ws.Range("E22").Value = "  +8.16%  "
$ws.Range("D23").Value = "13.23"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "27.370.26"
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("D25").Value = "2.470"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "2.664"
$ws.Range("E26").Value = "  +8.61%  "
$ws.Range("D27").Value = "21.94"
$ws.Range("E27").Value = "  +14.11%  "
$ws.Range("D28").Value = "1.470"
$ws.Range("E28").Value = "  +9.71%  "
$ws.Range("D29").Value = "154.47"
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("D30").Value = "2.020.53"
$ws.Range("E30").Value = "  +6.02%  "
$ws.Range("D31").Value = "135.90"
$ws.Range("E31").Value = "  +5.01%  "
$ws.Range("D32").Value = "6.317"
$ws.Range("E32").Value = "  +6.47%  "
$ws.Range("D33").Value = "4.031"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").Value = "13.90"
$ws.Range("E34").Value = "  +7.99%  "
$ws.Range("D35").Value = "0.08739"
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("D36").Value = "1.696"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "5.619"
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("D38").Value = "0.6969"
$ws.Range("E38").Value = "  +13.04%  "
$ws.Range("D39").Value = "0.2269"
$ws.Range("E39").Value = "  +6.11%  "
$ws.Range("D40").Value = "0.02405"
$ws.Range("E40").Value = "  +4.87%  "
$ws.Range("D41").Value = "0.06471"
$ws.Range("E41").Value = "  +3.74%  "
$ws.Range("D42").Value = "8.920"
$ws.Range("E42").Value = "  +3.89%  "
$ws.Range("D43").Value = "1.299"
$ws.Range("E43").Value = "  +5.69%  "
$ws.Range("D44").Value = "14.76"
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").Value = "0.6534"
$ws.Range("E45").Value = "  +10.83%  "
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "4.043"
$ws.Range("E47").Value = "  +5.29%  "
$ws.Range("D48").Value = "2.175"
$ws.Range("E48").Value = "  +7.69%  "
$ws.Range("D49").Value = "132.55"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("D50").Value = "0.07330"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").Value = "80.40"
$ws.Range("E51").Value = "  +4.34%  "
